# Adds a new "Sheet3" worksheet after the existing sheets, gives it a
# plain numeric value in A1, and updates Sheet2!A1's formula so that it
# references the new sheet instead of the literal (10/5) term. Also
# restores the original active sheet/selection afterwards.

$wb = $excel.ActiveWorkbook

# Remember which sheet was active before we start touching things.
$originalActiveSheetName = $wb.ActiveSheet.Name

# --- Add Sheet3 as the last worksheet -------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"
$ws3.Range("A1").Value = 2

# --- Update Sheet2!A1 formula to reference Sheet3!A1 -----------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Formula = "=POWER(2, 5)+8+Sheet3!A1"

# Record the new active cell (D39) on Sheet2.
$ws2.Activate()
$ws2.Range("D39").Select()

# --- Restore the workbook's original active sheet --------------------------
$wb.Worksheets.Item($originalActiveSheetName).Activate()
